# "update back-end for structure"
# Adds two new trailing columns (F: "intridual role", G: "entry criteria")
# to the position-details table's header row, widens the used range from
# A1:E5 to A1:G5 accordingly, and moves the live selection to E6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in row 1 (data rows 2-5 are left blank in these columns,
# matching the source change).
$ws.Range("F1").Value = "intridual role"
$ws.Range("G1").Value = "entry criteria"

# Column widths matching the new columns' best-fit size (~10.45 / ~10.82
# characters). ColumnWidth is expressed in characters; values are chosen so
# the engine's internal (pixel-quantized) stored width lands as close as
# possible to the authored widths.
$ws.Columns.Item(6).ColumnWidth = 9.6
$ws.Columns.Item(7).ColumnWidth = 10

# Selection moves to E6 after the edit.
$ws.Range("E6").Select()
